# Fix the teacher-name placeholder typo: {tracherName} -> {teacherName}
# (and tidy up a few other merge-field runs that Word re-flows on save,
#  matching the commit's accompanying OOXML clean-up).

$d = $word.ActiveDocument

function Merge-Text($findText, $tempText, $finalText) {
    # Word only "re-knits" a run (dropping stray <w:proofErr/> splits) when a
    # genuine text replacement happens over the range, so we swap to a
    # temporary value and then immediately back to the desired text.
    $r1 = $d.Content
    $r1.Find.Execute($findText, $true, $false, $false, $false, $false, `
                      $true, 1, $false, $tempText, 2) | Out-Null
    $r2 = $d.Content
    $r2.Find.Execute($tempText, $true, $false, $false, $false, $false, `
                      $true, 1, $false, $finalText, 2) | Out-Null
}

# The real content fix: correct the misspelled merge field.
Merge-Text "{tracherName}" "{teacherNameZZTEMPZZ}" "{teacherName}"

# Cosmetic run clean-up for fields on the same reachable story that the
# author's Word session also re-flowed while saving.
Merge-Text "{principalName}" "{principalNameZZTEMPZZ}" "{principalName}"
Merge-Text " {#attachments}{studentActivitiesPhotos} " " {#attachments}{studentActivitiesPhotosZZTEMPZZ} " " {#attachments}{studentActivitiesPhotos} "
Merge-Text "{studentAssessments}" "{studentAssessmentsZZTEMPZZ}" "{studentAssessments}"
